# Weekly refresh of the "Espinaca" price series:
#   - A brand-new observation is inserted as the new row 10
#     (date 2022-06-28 / 44740, volumen 35, precios 14000, $/Kg 1400).
#   - Every existing record that used to live in rows 10-25 is pushed
#     down by one row (10->11, 11->12, ... 25->26), so the table grows
#     from A1:R25 to A1:R26.
#
# We walk the affected rows bottom-up (25 -> 10) and copy each whole
# row onto the row below it so that sources are never overwritten
# before they are read. Column D (Fecha) also needs its date number
# format carried along since the destination row didn't have it yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstShiftRow = 10
$lastShiftRow  = 25
$lastCol       = 18   # column R

for ($row = $lastShiftRow; $row -ge $firstShiftRow; $row--) {
    for ($col = 1; $col -le $lastCol; $col++) {
        $src = $ws.Cells.Item($row, $col)
        $dst = $ws.Cells.Item($row + 1, $col)

        if ($col -eq 4) {
            $dst.NumberFormat = $src.NumberFormat
        }

        $dst.Value2 = $src.Value2
    }
}

# New observation for row 10 (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
$ws.Cells.Item(10, 4).Value2  = 44740   # Fecha
$ws.Cells.Item(10, 10).Value2 = 35      # Volumen
$ws.Cells.Item(10, 11).Value2 = 14000   # Precio minimo
$ws.Cells.Item(10, 12).Value2 = 14000   # Precio maximo
$ws.Cells.Item(10, 13).Value2 = 14000   # Precio promedio ponderado
$ws.Cells.Item(10, 16).Value2 = 1400    # Precio $/Kg
